# Adds a new "2022-Q4" sheet (cloned from "2022-Q3") with fresh quarterly
# fund data, and updates the "总计" (summary) sheet with the new quarter's
# row, shifting the older rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by duplicating "2022-Q3" (position 2)
#    so it inherits the same layout/styles, then insert it right before
#    the old sheet (i.e. it becomes the new position 2).
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

function Set-TextValue($ws, $cellRef, $text) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $text
    $ws.Range($cellRef).Style = "Normal"
}

# Row 2: 513090 / 易方达中证香港证券投资主题ETF
Set-TextValue $q4 "D2" "11.28"
Set-TextValue $q4 "E2" "97.07"
Set-TextValue $q4 "F2" "6.72"
Set-TextValue $q4 "G2" "0.7580"
$q4.Range("H2").Value = 8

# Row 3: 004497 / 前海开源多元策略灵活配置混合C
Set-TextValue $q4 "D3" "1.79"
Set-TextValue $q4 "E3" "79.66"
Set-TextValue $q4 "F3" "4.07"
Set-TextValue $q4 "G3" "0.0729"
$q4.Range("H3").Value = 6

# Row 4: 004496 / 前海开源多元策略灵活配置混合A
Set-TextValue $q4 "D4" "1.30"
Set-TextValue $q4 "E4" "79.66"
Set-TextValue $q4 "F4" "4.07"
Set-TextValue $q4 "G4" "0.0529"
$q4.Range("H4").Value = 6

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: add the 2022-Q4 row at the top of
#    the data (row 2) and push the rest down by one row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$labels = @("2022-Q4","2022-Q3","2022-Q2","2022-Q1","2021-Q4","2021-Q3","2021-Q2","2021-Q1","2020-Q4")
$values = @(0.88,0.88,1.04,1.15,1.12,1.08,0.82,0.55,0.5600000000000001)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $total.Range("A$row").Value = $i
    $total.Range("A$row").Style = "Normal"
    $total.Cells.Item($row, 1).Style = $total.Range("A1").Style
    $total.Range("B$row").Value = $labels[$i]
    $total.Range("C$row").Value = 3
    $total.Range("D$row").Value = $values[$i]
}

# ---------------------------------------------------------------------
# 3) Keep the originally-active sheet ("2020-Q4", always the last tab)
#    selected, matching the source workbook's tab state.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
